$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.160.70'
$ws.Range('E2').Value = '  -4.28%  '
$ws.Range('D3').Value = '2.223.70'
$ws.Range('E3').Value = '  -5.56%  '
$ws.Range('E4').Value = '  -0.15%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '318.59'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -3.70%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '98.69'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -7.44%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.580'
$c.Style = "Normal"
$ws.Range('E7').Value = '  -8.19%  '
$ws.Range('E8').Value = '  -0.16%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.565'
$c.Style = "Normal"
$ws.Range('E9').Value = '  -7.69%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '36.90'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -9.27%  '
$ws.Range('E11').Value = '  -3.38%  '
$ws.Range('E12').Value = '  -9.75%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '7.67'
$c.Style = "Normal"
$ws.Range('E13').Value = '  -8.62%  '
$ws.Range('E14').Value = '  -1.87%  '
$ws.Range('D15').Value = '2.562.91'
$ws.Range('E15').Value = '  -5.74%  '
$ws.Range('E16').Value = '  -11.28%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '14.33'
$c.Style = "Normal"
$ws.Range('E17').Value = '  -6.23%  '
$ws.Range('D18').Value = '2.227.44'
$ws.Range('E18').Value = '  -4.94%  '
$ws.Range('D19').Value = '43.064.89'
$ws.Range('E19').Value = '  -4.51%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '13.71'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -9.45%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '6.56'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -9.56%  '
$ws.Range('B22').Value = 'ShibaInu'
$ws.Range('C22').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D22').Value = '0.0₃0962'
$ws.Range('E22').Value = '  -9.16%  '
$ws.Range('E23').Value = '  -11.76%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '65.14'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -10.54%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '236.36'
$c.Style = "Normal"
$ws.Range('E25').Value = '  -8.45%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '2.20'
$c.Style = "Normal"
$ws.Range('E26').Value = '  -3.58%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -0.03%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '4.05'
$c.Style = "Normal"
$ws.Range('E28').Value = '  +1.53%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '10.06'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -10.60%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '2.21'
$c.Style = "Normal"
$ws.Range('E30').Value = '  -3.51%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '36.89'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +0.68%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '6.39'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -13.38%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '20.23'
$c.Style = "Normal"
$ws.Range('E33').Value = '  -8.47%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.0864'
$c.Style = "Normal"
$ws.Range('E34').Value = '  -9.90%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '157.80'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -5.38%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '3.32'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +1.35%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '2.65'
$c.Style = "Normal"
$ws.Range('E37').Value = '  -6.15%  '
$ws.Range('E38').Value = '  -8.17%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '1.85'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -3.89%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '4.39'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -6.65%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.104'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -10.24%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '3.68'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -7.01%  '
$ws.Range('E43').Value = '  -9.47%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '14.20'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +10.59%  '
$ws.Range('E45').Value = '  -0.17%  '
$ws.Range('D46').Value = '1.749.10'
$ws.Range('E46').Value = '  -6.90%  '
$ws.Range('E47').Value = '  -10.84%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '8.96'
$c.Style = "Normal"
$ws.Range('E48').Value = '  -3.35%  '
$ws.Range('B49').Value = 'BitcoinSV'
$ws.Range('C49').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '83.33'
$c.Style = "Normal"
$ws.Range('E49').Value = '  -12.90%  '
$ws.Range('E50').Value = '  -13.27%  '
$ws.Range('E51').Value = '  -12.94%  '
